$wb = $excel.ActiveWorkbook

# The same update applies to two worksheets: "展览" (sheet1) and "全部类型" (sheet4).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1230
    $ws.Range("F4").Value = 1453
    $ws.Range("F6").Value = 6118
}
